$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2772
$ws.Range("J18").Value = 6316
$ws.Range("L18").Value = 6316
$ws.Range("N18").Value = -6884
$ws.Range("H29").Value = 2203.8
$ws.Range("J29").Value = 3373
$ws.Range("L29").Value = 10119
$ws.Range("N29").Value = -10681
$ws.Range("H38").Value = 7209.0713
$ws.Range("I38").Value = 7209.0713
$ws.Range("K38").Value = 21627.2139
$ws.Range("M38").Value = -21255.2139
$ws.Range("H43").Value = 1000
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -931
$ws.Range("N43").ClearContents()
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H129").Value = 2100
$ws.Range("J129").Value = 15000
$ws.Range("L129").Value = 45000
$ws.Range("N129").Value = -55000
$ws.Range("H132").Value = 5570.9
$ws.Range("I132").Value = 5967.778
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 17903.334
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -15373.334
$ws.Range("N132").Value = -11057
$ws.Range("H135").Value = 1229.3
$ws.Range("I135").Value = 1047
$ws.Range("J135").Value = 1654.6666
$ws.Range("K135").Value = 9423
$ws.Range("L135").Value = 14891.9994
$ws.Range("M135").Value = -6888
$ws.Range("N135").Value = -19961.9994
$ws.Range("H138").Value = 2795.4546
$ws.Range("J138").Value = 3000
$ws.Range("L138").Value = 9000
$ws.Range("N138").Value = -19280
$ws.Range("H141").Value = 724.2857
$ws.Range("I141").Value = 661.6667
$ws.Range("K141").Value = 1985.0001
$ws.Range("M141").Value = 3194.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 500
$ws.Range("I38").Value = 500
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 500
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -33
$ws.Range("N38").ClearContents()
$ws.Range("H45").Value = 1434.2222
$ws.Range("I45").Value = 651.3333
$ws.Range("K45").Value = 651.3333
$ws.Range("M45").Value = -274.3333
$ws.Range("H63").Value = 2500
$ws.Range("J63").Value = 2500
$ws.Range("L63").Value = 2500
$ws.Range("N63").Value = -3872
$ws.Range("H66").Value = 2500
$ws.Range("J66").Value = 2500
$ws.Range("L66").Value = 12500
$ws.Range("N66").Value = -19364

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 7459.6665
$ws.Range("I37").Value = 350
$ws.Range("J37").Value = 11014.5
$ws.Range("K37").Value = 350
$ws.Range("L37").Value = 11014.5
$ws.Range("M37").Value = -213
$ws.Range("N37").Value = -11288.5
$ws.Range("H80").Value = 1934.3846
$ws.Range("I80").Value = 1739
$ws.Range("K80").Value = 1739
$ws.Range("M80").Value = -741
$ws.Range("H83").Value = 1934.3846
$ws.Range("I83").Value = 1739
$ws.Range("K83").Value = 8695
$ws.Range("M83").Value = -3703

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 2316.6667
$ws.Range("I32").Value = 2316.6667
$ws.Range("K32").Value = 2316.6667
$ws.Range("M32").Value = -2000.6667
$ws.Range("H41").Value = 8806.429
$ws.Range("I41").Value = 1941
$ws.Range("J41").Value = 49999
$ws.Range("K41").Value = 1941
$ws.Range("L41").Value = 49999
$ws.Range("M41").Value = -1513
$ws.Range("N41").Value = -50855
$ws.Range("H99").Value = 11714.6
$ws.Range("I99").Value = 11714.6
$ws.Range("K99").Value = 11714.6
$ws.Range("M99").Value = -10216.6
$ws.Range("H126").Value = 11714.6
$ws.Range("I126").Value = 11714.6
$ws.Range("K126").Value = 35143.8
$ws.Range("M126").Value = -32673.8
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1001
$ws.Range("I16").Value = 1001
$ws.Range("K16").Value = 3003
$ws.Range("M16").Value = -2830
$ws.Range("H23").Value = 679.6667
$ws.Range("I23").Value = 519.5
$ws.Range("J23").Value = 1000
$ws.Range("K23").Value = 1558.5
$ws.Range("L23").Value = 3000
$ws.Range("M23").Value = -1323.5
$ws.Range("N23").Value = -3470

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 21008
$ws.Range("I22").Value = 21008
$ws.Range("K22").Value = 21008
$ws.Range("M22").Value = -20479
$ws.Range("H24").Value = 2100
$ws.Range("I24").Value = 2100
$ws.Range("K24").Value = 2100
$ws.Range("M24").Value = -1927
$ws.Range("H46").Value = 3030.625
$ws.Range("I46").Value = 1314.1428
$ws.Range("K46").Value = 1314.1428
$ws.Range("M46").Value = -1158.1428

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 408.63635
$ws.Range("I22").Value = 377.22223
$ws.Range("K22").Value = 377.22223
$ws.Range("M22").Value = -82.22223000000002
$ws.Range("H27").Value = 408.63635
$ws.Range("I27").Value = 377.22223
$ws.Range("K27").Value = 377.22223
$ws.Range("M27").Value = -270.22223
$ws.Range("H31").Value = 26500
$ws.Range("I31").Value = 5000
$ws.Range("J31").Value = 48000
$ws.Range("K31").Value = 5000
$ws.Range("L31").Value = 48000
$ws.Range("M31").Value = -4752
$ws.Range("N31").Value = -48496
$ws.Range("H55").Value = 396.75
$ws.Range("I55").Value = 396.75
$ws.Range("K55").Value = 396.75
$ws.Range("M55").Value = -223.75
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 10642.857
$ws.Range("I136").Value = 10642.857
$ws.Range("K136").Value = 31928.571
$ws.Range("M136").Value = -29378.571
